# Remove "scATACseq" from the assay_type and transposition_method lookup
# lists (and thereby from the shared-strings table), and fix up the
# dataValidation rules on the "Export as TSV" sheet that reference those
# lists.

$wb = $excel.ActiveWorkbook

# --- assay_type list: SNARE-seq2 / scATACseq / sciATACseq / snATACseq ---
# Remove the "scATACseq" row (row 2); sciATACseq/snATACseq shift up.
$wsAssayType = $wb.Worksheets.Item("assay_type list")
$wsAssayType.Range("A2").EntireRow.Delete()

# --- transposition_method list: SNARE-Seq2-AC / scATACseq / bulkATACseq / snATACseq / sciATACseq ---
# Remove the "scATACseq" row (row 2); the rest shift up.
$wsTranspositionMethod = $wb.Worksheets.Item("transposition_method list")
$wsTranspositionMethod.Range("A2").EntireRow.Delete()

# --- Fix up the data validations on "Export as TSV" that reference the
#     two lists above, since their row counts shrank by one each. ---
$wsTsv = $wb.Worksheets.Item("Export as TSV")

$assayTypeValidation = $wsTsv.Range("L2:L1048576").Validation
$assayTypeValidation.Formula1 = "'assay_type list'!`$A`$1:`$A`$3"
$assayTypeValidation.ErrorMessage = "Value must be one of: SNARE-seq2 / sciATACseq / snATACseq."

$transpositionMethodValidation = $wsTsv.Range("Z2:Z1048576").Validation
$transpositionMethodValidation.Formula1 = "'transposition_method list'!`$A`$1:`$A`$4"
$transpositionMethodValidation.ErrorMessage = "Value must be one of: SNARE-Seq2-AC / bulkATACseq / snATACseq / sciATACseq."
